$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row for "RM 232" (row 26) entirely; rows below shift up.
$ws.Rows.Item(26).EntireRow.Delete()

# After the above deletion, "SC 92" (originally row 28) is now at row 27.
# Delete it entirely too; rows below shift up again.
$ws.Rows.Item(27).EntireRow.Delete()

# Cell-level value edits within the remaining data.
$ws.Range("D19").Value = -15.5
$ws.Range("D21").ClearContents()
$ws.Range("D23").Value = -13.9

# "SC 101" row (now row 27) loses its D value.
$ws.Range("D27").ClearContents()

# "SC 232" row (now row 33) gains a D value.
$ws.Range("D33").Value = -14.1
